$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet LOGT2 -> logt2
$ws.Name = "logt2"

# Try to set the tab split ratio (bookViews workbookView tabRatio 456 -> 400)
$win = $excel.ActiveWindow
$win.TabRatio = 0.4

# Date in F1 changes from 41905 to 41909
$ws.Range("F1").Value = 41909

# Insert a new blank row above row 5 -- shifts header row 5->6, data rows 6->7 and 7->8,
# and the trailing summary row 9->10 (formula range auto-adjusts to SUM(E7:E8)).
$ws.Rows.Item(5).Insert()

# The old F6 (now F7) task-count value of 1 is cleared (blank in new layout).
$ws.Range("F7").ClearContents()

# The old F7 (now F8) task-count value changes from 6 to 12.
$ws.Range("F8").Value = 12

# Append the two new log rows (9 and 10).
$ws.Range("A9").Value = 41912
$ws.Range("B9").Value = 0.96875
$ws.Range("C9").Value = 0.993055555555555
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 35
$ws.Range("H9").Value = "Tuvimos una reunión de equipo para la asignación de las tareas del ciclo #1."

$ws.Range("A10").Value = 41913
$ws.Range("B10").Value = 0.833333333333333
$ws.Range("C10").Value = 0.930555555555555
$ws.Range("D10").Value = 0
# Row 10's E cell previously held the old trailing summary formula (shifted down by the
# insert above); overwrite it with the real per-row delta-time formula for the new entry.
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("F10").Value = "13,14,15"
$ws.Range("H10").Value = "Participé en el diagrama de caso de uso y escenarios de atributos de calidad"

# The new summary formula lives in the freshly inserted row 5, covering the full data range.
$ws.Range("E5").Formula = "=SUM(E7:E10)/60"

# Newly typed-in cells (F7/F8/F9/F10 "Assembly" counts) are right-aligned, and the new
# blank G9/G10 "Assembly" cells are center-aligned -- matching the sheet's new column
# formatting for those positions.
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("F10").HorizontalAlignment = -4152
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G10").HorizontalAlignment = -4108

# Column width tweaks (the sheet was narrowed slightly across the board).
$ws.Columns.Item(1).ColumnWidth = 11.083333333372138
$ws.Range("B1:C1").ColumnWidth = 11.083333333372138
$ws.Columns.Item(4).ColumnWidth = 15.916666666671517
$ws.Columns.Item(5).ColumnWidth = 10.08333333338669
$ws.Range("F1:G1").ColumnWidth = 11.083333333372138
$ws.Columns.Item(8).ColumnWidth = 46.41666666670062
